$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.113767862319946
$ws.Range("B1").Value = 2.075846433639526
$ws.Range("C1").Value = 1.973269581794739
$ws.Range("D1").Value = 2.629579782485962
$ws.Range("E1").Value = 5.166424751281738
